$wb = $excel.ActiveWorkbook

# =====================================================================
# 1. Add the new worksheet "Patientputvital" after "PatientPut"
# =====================================================================
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "Patientputvital"

$wsAdmin = $wb.Worksheets.Item("AdminLogin")
$wsPost  = $wb.Worksheets.Item("PatientPost")
$wsPut   = $wb.Worksheets.Item("PatientPut")
$wsVital = $wb.Worksheets.Item("Patientputvital")

# =====================================================================
# 2. AdminLogin sheet - add row 3 + column E test data
# =====================================================================
$wsAdmin.Range("A3").Value = "test"
$wsAdmin.Range("B3").Value = "team123@gmail.com"
$wsAdmin.Hyperlinks.Add($wsAdmin.Range("B3"), "mailto:team123@gmail.com") | Out-Null
$wsAdmin.Range("B2").Copy() | Out-Null
$wsAdmin.Range("B3").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$wsAdmin.Range("E2").Value = "Valid data"
$wsAdmin.Range("E3").Value = "invalid data"

# =====================================================================
# 3. PatientPost sheet
# =====================================================================
# 3a. Fix the date-format styling on column H (numFmt 165 -> numFmt 166)
$wsPost.Range("H1").NumberFormat = "yyyy\-mm\-dd;@"
$wsPost.Range("H3").NumberFormat = "yyyy\-mm\-dd;@"
$wsPost.Range("H2").NumberFormat = "yyyy\-mm\-dd;@"
$wsPost.Range("H2").Font.Name = "Consolas"
$wsPost.Range("H2").Font.Size = 10
$wsPost.Range("H2").Font.Color = 0

# 3b. New rows 6 and 7 (negative-test vitals data)
$wsPost.Range("A6").Value = "Sue12"
$wsPost.Range("B6").Value = "katie1"
$wsPost.Range("C6").Value = 3457899
$wsPost.Range("D6").Value = "jaf54@gmail.com"
$wsPost.Hyperlinks.Add($wsPost.Range("D6"), "mailto:jaf54@gmail.com") | Out-Null
$wsPost.Range("E6").Value = "SOY"
$wsPost.Range("F6").Value = "Ve%gan"
$wsPost.Range("G6").Value = "Indian"
$wsPost.Range("H6").NumberFormat = "@"
$wsPost.Range("H6").Value = "123-12-30"
$wsPost.Range("H6").NumberFormat = "yyyy\-mm\-dd;@"
$wsPost.Range("K6").Value = "invalid Mandatory data"

$wsPost.Range("A7").Value = "Ram"
$wsPost.Range("B7").Value = "Swamy"
$wsPost.Range("C7").Value = 8408305647
$wsPost.Range("D7").Value = "Ramy@gmail.com"
$wsPost.Hyperlinks.Add($wsPost.Range("D7"), "mailto:Ramy@gmail.com") | Out-Null
$wsPost.Range("E7").Value = "SOY"
$wsPost.Range("F7").Value = "Vega"
$wsPost.Range("G7").Value = "Tamil"
$wsPost.Range("H7").NumberFormat = "yyyy\-mm\-dd;@"
$wsPost.Range("H7").Value = 31865
$wsPost.Range("K7").Value = "Valid Mandatory invalid additional field"

$wsPost.Range("G10").Select() | Out-Null

# =====================================================================
# 4. PatientPut sheet
# =====================================================================
# 4a. Duplicate row 2 pattern into rows 3 and 4, overwriting old content
$wsPut.Range("A2:H2").Copy() | Out-Null
$wsPut.Range("A3").PasteSpecial(-4104) | Out-Null
$wsPut.Range("A2:H2").Copy() | Out-Null
$wsPut.Range("A4").PasteSpecial(-4104) | Out-Null
$wsPut.Range("A2:H2").Copy() | Out-Null
$wsPut.Range("A5").PasteSpecial(-4104) | Out-Null
$excel.CutCopyMode = $false

# 4b. Row-height / dyDescent reset for rows 3 & 4 (remove the 17.25 custom height)
$wsPut.Rows("3").RowHeight = 15
$wsPut.Rows("4").RowHeight = 15

# 4c. Row 5 differs in E/F/G (and E5 gets a hyperlink instead of plain text)
$wsPut.Range("E5").Value = "Wal@nut"
$wsPut.Hyperlinks.Add($wsPut.Range("E5"), "mailto:Wal@nut") | Out-Null
$wsPut.Range("F5").Value = "West Vegan"
$wsPut.Range("G5").Value = "Indiaaa"

# 4d. Row 6 - blank styled row (copy formats only from row 2)
$wsPut.Range("C2:H2").Copy() | Out-Null
$wsPut.Range("C6").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$wsPut.Range("C6:H6").ClearContents()

# 4e. K / L columns
$wsPut.Range("K2").Value = "Positive"
$wsPut.Range("L2").Value = "update any data with existing pdf file (update either contact number or food preference)"

$wsPut.Range("K3").Value = "Positive"
$wsPut.Range("L3").Value = "update any data focus only Mandatory and additional details  without adding pdf file  (update either email or DOB) "

$wsPut.Range("K4").Value = "Positive"
$wsPut.Range("L4").Value = "Updating only pdf file ,no need to change the data(NO data update)"

$wsPut.Range("K5").Value = "Negative"
$wsPut.Range("L5").Value = "invalid additional field"

$wsPut.Range("G8").Select() | Out-Null

# =====================================================================
# 5. Patientputvital sheet - formats only, copied from PatientPut row 1/3
# =====================================================================
$wsPut.Range("H1").Copy() | Out-Null
$wsVital.Range("H1").PasteSpecial(-4122) | Out-Null

$wsPut.Range("D3").Copy() | Out-Null
$wsVital.Range("D2").PasteSpecial(-4122) | Out-Null

$wsPut.Range("F3").Copy() | Out-Null
$wsVital.Range("F2").PasteSpecial(-4122) | Out-Null

$wsPut.Range("G3").Copy() | Out-Null
$wsVital.Range("G2").PasteSpecial(-4122) | Out-Null

$wsPut.Range("H3").Copy() | Out-Null
$wsVital.Range("H2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$wsVital.Range("F32").Select() | Out-Null

# =====================================================================
# 6. Final selections / active sheet
# =====================================================================
$wsAdmin.Activate()
$wsAdmin.Range("E7").Select() | Out-Null
